# Binary Search: Kth smallest price
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")
$ws.Activate()

# New row of data
$ws.Range("B9").Value = 7
$ws.Range("D9").Value = "Kth Smallest Price"

$cell = $ws.Range("E9")
$cell.Value = "KthPrice - Problem | Scaler Academy"
$ws.Hyperlinks.Add(
    $cell,
    "https://www.scaler.com/academy/mentee-dashboard/class/30365/homework/problems/872?navref=cl_tt_nv",
    [Type]::Missing,
    [Type]::Missing,
    "KthPrice - Problem | Scaler Academy"
) | Out-Null

# Match formatting of the row above (row height + styles)
$ws.Rows.Item(9).RowHeight = 28.8
$ws.Range("B9").Style = $ws.Range("B8").Style
$ws.Range("D9").Style = $ws.Range("D8").Style
$ws.Range("E9").Style = $ws.Range("E8").Style

$ws.Range("F9").Select()
